{"js": "// Remove the \"Appendix: Quick prototype\" section (heading, the intro\n// line, the \"Figure: PDF page 1\" caption and the embedded screenshot\n// paragraph) that used to sit between the first \"Appendix: Links\" text\n// paragraph and the real \"Appendix: Links\" Heading 2 section.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Appendix: Quick prototype\" Heading 2 paragraph - the start\n// of the block to remove.\nlet startIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].style === \"Heading 2\" && items[i].text === \"Appendix: Quick prototype\") {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex !== -1) {\n  // Locate the next \"Appendix: Links\" Heading 2 paragraph after the start\n  // - the paragraph right after the block we need to delete (exclusive).\n  let endIndex = -1;\n  for (let j = startIndex + 1; j < items.length; j++) {\n    if (items[j].style === \"Heading 2\" && items[j].text === \"Appendix: Links\") {\n      endIndex = j;\n      break;\n    }\n  }\n\n  if (endIndex !== -1) {\n    // Delete from the end backwards so earlier indices stay valid.\n    for (let k = endIndex - 1; k >= startIndex; k--) {\n      items[k].delete();\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the \"Appendix: Quick prototype\" section - the Heading 2 title,\n# the blank line, the \"Figure: PDF page 1\" caption and the paragraph\n# holding the embedded screenshot - leaving the real \"Appendix: Links\"\n# Heading 2 section (and everything before the prototype section)\n# untouched.\n$d = $word.ActiveDocument\n\n# Locate the \"Appendix: Quick prototype\" heading - the start of the block\n# to remove.\n$startSearch = $d.Content\n$startSearch.Find.ClearFormatting()\n$startFound = $startSearch.Find.Execute(\"Appendix: Quick prototype\")\n\nif ($startFound) {\n    $startPara = $startSearch.Paragraphs.Item(1)\n    $startRange = $d.Range($startPara.Range.Start, $startPara.Range.Start)\n\n    # From just past the heading, find the next \"Appendix: Links\" heading\n    # - the paragraph that should immediately follow once the block is\n    # gone.\n    $endSearch = $d.Range($startSearch.End, $d.Content.End)\n    $endSearch.Find.ClearFormatting()\n    $endFound = $endSearch.Find.Execute(\"Appendix: Links\")\n\n    if ($endFound) {\n        $endPara = $endSearch.Paragraphs.Item(1)\n        $endRange = $d.Range($endPara.Range.Start, $endPara.Range.Start)\n\n        # Delete everything from the start of the prototype heading up to\n        # (but not including) the start of the following \"Appendix: Links\"\n        # heading.\n        $delRange = $d.Range($startRange.Start, $endRange.Start)\n        $delRange.Delete()\n    }\n}\n"}
